$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows right before the old row 21 (so the new
# work-log entry lands at row 21, the "Sum" row re-targets C3:C24, and
# everything below shifts down by two rows, matching the target diff).
$ws.Rows("21:22").Insert()

# The freshly inserted rows 21/22 don't inherit the surrounding blank-row
# formatting automatically, so copy it over from rows 23/24 (which *are*
# the old, correctly-styled rows 21/22, now shifted down).
$ws.Range("A23:E23").Copy()
[void]$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A24:E24").Copy()
[void]$ws.Range("A22:E22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New journal entry for the UiController work.
$ws.Cells.Item(21, 1).Value = "DataMapper, Validator, InverseStringResponse, InverseStringRequest, тесты"
$ws.Cells.Item(21, 2).Value = 45684
$ws.Cells.Item(21, 3).Value = 0.125
$ws.Cells.Item(21, 4).Value = "UiController, InverseServiceImplTest"

# Cosmetic view updates that accompanied the edit.
[void]$ws.Range("D22").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1

# Column D widened slightly to fit the new, longer component description.
$ws.Columns.Item(4).ColumnWidth = 34.416666666666664
